$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) held values one day off (source stats used a
# different day boundary than the rest of the season folder naming).
# Replace "7-1-2011-12" with the corrected "2012-07-01" for every data
# row (BF2:BF31), keeping the cells as plain text (not auto-converted
# to a date serial) and leaving their formatting untouched.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.Formula = "=""2012-07-01"""
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0
